$wb = $excel.ActiveWorkbook

# --- Fix the "fraction" row formula on the mads_tightened sheet ---
# It was hard-coded to LOG10(0.6) instead of referencing the base value
# on the summary sheet (summary!B5), unlike every other row in this
# sheet (and unlike the equivalent row on the "mads" sheet).
$wsTightened = $wb.Worksheets.Item("mads_tightened")
$wsTightened.Range("B5").Formula = "=LOG10(summary!B5)"

# --- Restore the selected cell on each sheet ---
$wsSummary = $wb.Worksheets.Item("summary")
$wsSummary.Range("B6").Select()

$wsMads = $wb.Worksheets.Item("mads")
$wsMads.Range("B4").Select()

# Leave mads_tightened as the active sheet/selection, matching the
# original workbook state (it was the tabSelected sheet before and
# after the edit).
$wsTightened.Range("C16").Select()
